# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 1681333.4
$ws.Range("I54").Value = 2505000
$ws.Range("J54").Value = 34000
$ws.Range("K54").Value = 2505000
$ws.Range("L54").Value = 34000
$ws.Range("M54").Value = -2504514
$ws.Range("N54").Value = -34972

$ws.Range("H62").Value = 14126.25
$ws.Range("I62").Value = 8835
$ws.Range("J62").Value = 30000
$ws.Range("K62").Value = 8835
$ws.Range("L62").Value = 30000
$ws.Range("M62").Value = -8211
$ws.Range("N62").Value = -31248

$ws.Range("H65").Value = 14126.25
$ws.Range("I65").Value = 8835
$ws.Range("J65").Value = 30000
$ws.Range("K65").Value = 44175
$ws.Range("L65").Value = 150000
$ws.Range("M65").Value = -41055
$ws.Range("N65").Value = -156240

$ws.Range("H125").Value = 12457512
$ws.Range("I125").Value = 654.3333
$ws.Range("J125").Value = 18685942
$ws.Range("K125").Value = 5888.9997
$ws.Range("L125").Value = 168173478
$ws.Range("M125").Value = -3428.9997
$ws.Range("N125").Value = -168178398

$ws.Range("H135").Value = 8799.571
$ws.Range("I135").Value = 10016.167
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 90145.503
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -87610.503
$ws.Range("N135").Value = -18570

$ws.Range("H138").Value = 324941.56
$ws.Range("I138").Value = 1467210.5
$ws.Range("J138").Value = 2763.141
$ws.Range("K138").Value = 4401631.5
$ws.Range("L138").Value = 8289.423000000001
$ws.Range("M138").Value = -4396491.5
$ws.Range("N138").Value = -18569.423

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17224.262
$ws.Range("I32").Value = 2008.08
$ws.Range("J32").Value = 57266.844
$ws.Range("K32").Value = 2008.08
$ws.Range("L32").Value = 57266.844
$ws.Range("M32").Value = -1721.08
$ws.Range("N32").Value = -57840.844

$ws.Range("H74").Value = 9902.4
$ws.Range("I74").Value = 2370.6667
$ws.Range("J74").Value = 21200
$ws.Range("K74").Value = 2370.6667
$ws.Range("L74").Value = 21200
$ws.Range("M74").Value = -1496.6667
$ws.Range("N74").Value = -22948

$ws.Range("H77").Value = 9902.4
$ws.Range("I77").Value = 2370.6667
$ws.Range("J77").Value = 21200
$ws.Range("K77").Value = 11853.3335
$ws.Range("L77").Value = 106000
$ws.Range("M77").Value = -7485.333500000001
$ws.Range("N77").Value = -114736

$ws.Range("H97").Value = 8747.5
$ws.Range("I97").Value = 10372
$ws.Range("J97").Value = 625
$ws.Range("K97").Value = 10372
$ws.Range("L97").Value = 625
$ws.Range("M97").Value = -9876

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3934.7856
$ws.Range("I134").Value = 2911
$ws.Range("J134").Value = 6494.25
$ws.Range("K134").Value = 8733
$ws.Range("L134").Value = 19482.75
$ws.Range("M134").Value = -6198

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2035.2941
$ws.Range("I58").Value = 1111.3
$ws.Range("J58").Value = 3355.2856
$ws.Range("K58").Value = 1111.3
$ws.Range("L58").Value = 3355.2856
$ws.Range("M58").Value = -908.3
$ws.Range("N58").Value = -3761.2856

$ws.Range("H134").Value = 2661.5
$ws.Range("I134").Value = 1447.45
$ws.Range("J134").Value = 5696.625
$ws.Range("K134").Value = 4342.35
$ws.Range("L134").Value = 17089.875
$ws.Range("M134").Value = -1807.35
$ws.Range("N134").Value = -22159.875

$ws.Range("H136").Value = 2035.2941
$ws.Range("I136").Value = 1111.3
$ws.Range("J136").Value = 3355.2856
$ws.Range("K136").Value = 3333.9
$ws.Range("L136").Value = 10065.8568
$ws.Range("M136").Value = -783.8999999999996
$ws.Range("N136").Value = -15165.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 946.6667
$ws.Range("I51").Value = 893.3333
$ws.Range("J51").Value = 1000
$ws.Range("K51").Value = 2679.9999
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = -2219.9999

$ws.Range("H113").Value = 16667514
$ws.Range("I113").Value = 389.6
$ws.Range("J113").Value = 20000940
$ws.Range("K113").Value = 1168.8
$ws.Range("L113").Value = 60002820
$ws.Range("M113").Value = 1001.2
$ws.Range("N113").Value = -60007160

$ws.Range("H122").Value = 669.1053000000001
$ws.Range("I122").Value = 302
$ws.Range("J122").Value = 838.53845
$ws.Range("K122").Value = 2718
$ws.Range("L122").Value = 7546.84605
$ws.Range("M122").Value = -268
$ws.Range("N122").Value = -12446.84605

$ws.Range("H131").Value = 1432.25
$ws.Range("I131").Value = 463.18182
$ws.Range("J131").Value = 1633.3773
$ws.Range("K131").Value = 1389.54546
$ws.Range("L131").Value = 4900.1319
$ws.Range("M131").Value = 3650.45454
$ws.Range("N131").Value = -14980.1319

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6313.5454
$ws.Range("I70").Value = 6494.3335
$ws.Range("J70").Value = 5500
$ws.Range("K70").Value = 6494.3335
$ws.Range("L70").Value = 5500
$ws.Range("M70").Value = -6224.3335

$ws.Range("H73").Value = 6313.5454
$ws.Range("I73").Value = 6494.3335
$ws.Range("J73").Value = 5500
$ws.Range("K73").Value = 6494.3335
$ws.Range("L73").Value = 5500
$ws.Range("M73").Value = -5558.3335

$ws.Range("H102").Value = 2450
$ws.Range("I102").Value = 1900
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 1900
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -278
$ws.Range("N102").Value = -6244

$ws.Range("H113").Value = 1530
$ws.Range("I113").Value = 1530
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1530
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 640
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 1236093.1
$ws.Range("I122").Value = 1390429.8
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 4171289.4
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = -4168839.4
$ws.Range("N122").Value = -9100

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4204.737
$ws.Range("I40").Value = 2347.5
$ws.Range("J40").Value = 4700
$ws.Range("K40").Value = 2347.5
$ws.Range("L40").Value = 4700
$ws.Range("M40").Value = -2211.5
$ws.Range("N40").Value = -4972

$ws.Range("H100").Value = 2933.2856
$ws.Range("I100").Value = 2679.8
$ws.Range("J100").Value = 3012.5
$ws.Range("K100").Value = 2679.8
$ws.Range("L100").Value = 3012.5
$ws.Range("M100").Value = -2138.8
$ws.Range("N100").Value = -4094.5

$ws.Range("H136").Value = 6870.6875
$ws.Range("I136").Value = 3054.3333
$ws.Range("J136").Value = 11777.429
$ws.Range("K136").Value = 9162.999899999999
$ws.Range("L136").Value = 35332.287
$ws.Range("M136").Value = -6612.999899999999
$ws.Range("N136").Value = -40432.287

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1395
$ws.Range("I81").Value = 1395
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 2790
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -1729
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 1395
$ws.Range("I84").Value = 1395
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 13950
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -8646
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 92453.91
$ws.Range("I122").Value = 112443.664
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 337330.992
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -334880.992
$ws.Range("N122").Value = -12400

$ws.Range("H136").Value = 15198361
$ws.Range("I136").Value = 18574688
$ws.Range("J136").Value = 4885.5
$ws.Range("K136").Value = 55724064
$ws.Range("L136").Value = 14656.5
$ws.Range("M136").Value = -55721514
$ws.Range("N136").Value = -19756.5
